$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, F, H, L, M, N across rows 2-25 (data rows 0-23)
$data = @{}

$data["B"] = @(2.462849570000458, 2.339734690346575, 2.26572530936636, 2.235961576800605, 2.231043157541535, 2.26532230645023, 2.420069291190828, 2.736227738671232, 2.976467311397016, 3.087535698457714, 3.129854200067484, 3.120728571409245, 3.091012052373628, 3.072843688059095, 2.96924485140255, 2.906148922517161, 2.870025376367607, 2.85782325214052, 2.912848227931818, 3.099733452099827, 3.223386689591791, 3.157251173789575, 2.909819001402354, 2.649316820617685)
$data["C"] = @(0.2267414144939721, 0.1971758706343962, 0.179039771927421, 0.1716527671450478, 0.1704263568152271, 0.1789401345724286, 0.2165431763653771, 0.2904566404548063, 0.3449235667931134, 0.369749989201182, 0.3791591112508854, 0.3771323240249558, 0.3705239203225688, 0.3664771351643026, 0.3433021610682658, 0.3290982369068161, 0.3209330907955632, 0.3181692734488308, 0.3306097893927245, 0.3724647477569647, 0.3998657928131593, 0.3852368157873229, 0.3299264138379385, 0.2704366529130198)
$data["D"] = @(0.1261810435928226, 0.1267737299606289, 0.1271821046592834, 0.1273596603933207, 0.1273898150171888, 0.1271844541887184, 0.126376148800496, 0.1251457441738282, 0.1244606681744855, 0.1241971384477267, 0.1241043155248676, 0.124123995780387, 0.1241893618618803, 0.124230309721014, 0.1244788630100473, 0.1246437018731399, 0.1247430384952253, 0.1247774479126562, 0.1246256857373851, 0.1241699726798231, 0.1239127892296992, 0.1240463160833656, 0.1246338166063481, 0.1254403567337263)
$data["E"] = @(0.05179948103848364, 0.05127343437879794, 0.05094355868143374, 0.05080738384047745, 0.05078466603803111, 0.05094172927922891, 0.05161951757817373, 0.05289489744467879, 0.05380036671946797, 0.0542057063608663, 0.05435827494959611, 0.05432545739713035, 0.05421827669375645, 0.05415250555848949, 0.05377374637718191, 0.05353972241713656, 0.0534044991452074, 0.05335860800103376, 0.05356469858211277, 0.05424978320900031, 0.0546921445419688, 0.0544565340000176, 0.05355340897676442, 0.05255550338582093)
$data["F"] = @(2.512590215698481, 2.463575717352057, 2.434976978526691, 2.423696730668482, 2.421846174276411, 2.43482333802568, 2.495378032804112, 2.626112249226622, 2.729641756673573, 2.778402473010175, 2.797109177730391, 2.793069543975179, 2.779936617551897, 2.771923928000433, 2.726488856198358, 2.699044174495356, 2.683415237091282, 2.67815034766241, 2.701949489261267, 2.783787484600197, 2.838685660933095, 2.809255337126018, 2.700635531387661, 2.589444330953796)
$data["H"] = @(0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429)
$data["L"] = @(0.1935801646048105, 0.1912509182809643, 0.1899219057063846, 0.189405717609489, 0.1893215376319404, 0.1899148414405403, 0.1927560197516343, 0.1991326153031281, 0.2043125403091608, 0.206777519720788, 0.2077266295704163, 0.2075215236528578, 0.2068552890541469, 0.2064492443253982, 0.2041536384624294, 0.2027732158908577, 0.2019894521868792, 0.2017258366205965, 0.202919106202728, 0.2070505523975612, 0.2098420896090261, 0.2083438108709572, 0.202853118486729, 0.1973209581740676)
$data["M"] = @(0.443544082336615, 0.4261271886266371, 0.4157147108890555, 0.4115419923478498, 0.4108533605601608, 0.4156581513565811, 0.437480078449255, 0.4825258030473236, 0.5170249657290498, 0.5330314114518657, 0.5391380859501709, 0.5378208790628918, 0.5335328992426085, 0.5309123106331768, 0.5159852336579434, 0.5069083224728317, 0.5017169206700842, 0.4999642373630024, 0.5078715295689733, 0.5347911477622915, 0.5526493957191789, 0.5430937461571972, 0.5074359794639349, 0.4700954177488583)
$data["N"] = @(2.342519244637273, 2.353551735120391, 2.361022628853064, 2.364241640070588, 2.364786677732539, 2.361065335703046, 2.346178194808672, 2.322545106758909, 2.308615053884523, 2.303032930292858, 2.301028456949766, 2.301455279218828, 2.302865824636356, 2.30374409051521, 2.308995116236858, 2.312410338694377, 2.314445672688137, 2.315146969109108, 2.312039429459304, 2.302448538720441, 2.296818138332554, 2.299764553755807, 2.312206893705721, 2.328338751546298)

foreach ($col in $data.Keys) {
    $colIndex = $ws.Range($col + "1").Column
    $values = $data[$col]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = 2 + $i
        $ws.Cells.Item($row, $colIndex).Value2 = $values[$i]
    }
}

Write-Host "Applied 380 kV case values"
